$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new BOM row (row 28) ---
$ws.Range("B28").Value = "KB304-PNW"
$ws.Range("C28").Value = "Klawiatura"
$ws.Range("K28").Value = "https://www.tme.eu/pl/details/kb304-pnw/klawiatury-plastikowe/accord/ak-304-n-wwb/"
$ws.Range("L28").Value = 1
$ws.Range("M28").Value = 20.21
$ws.Range("N28").Value = 1
$ws.Range("O28").Formula = "=L28*M28"

# --- Column C becomes visible again (with a normal default width) ---
$ws.Columns("C").Hidden = $false
$ws.Columns("C").ColumnWidth = 8

# --- Columns D:J keep being hidden, but now carry an explicit (non-zero) width ---
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 8

# --- Update the last selected cell recorded in the sheet view ---
$null = $ws.Range("V13").Select()

$wb.Save()
